# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '66.073.56'
Set-TextValue $ws.Range("E2") '  +0.06%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.174.01'
Set-TextValue $ws.Range("E3") '  -1.04%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '607.42'
Set-TextValue $ws.Range("E5") '  +0.91%  '

# Row 6
Set-TextValue $ws.Range("D6") '154.49'
Set-TextValue $ws.Range("E6") '  +0.45%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +0.04%  '

# Row 8
Set-TextValue $ws.Range("D8") '3.171.41'
Set-TextValue $ws.Range("E8") '  -1.02%  '

# Row 9
Set-TextValue $ws.Range("E9") '  +2.31%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.158'
Set-TextValue $ws.Range("E10") '  -1.18%  '

# Row 11
Set-TextValue $ws.Range("D11") '5.69'
Set-TextValue $ws.Range("E11") '  -6.79%  '

# Row 12
Set-TextValue $ws.Range("E12") '  +1.24%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.0000266'
Set-TextValue $ws.Range("E13") '  -2.09%  '

# Row 14
Set-TextValue $ws.Range("D14") '38.19'
Set-TextValue $ws.Range("E14") '  -3.21%  '

# Row 15
Set-TextValue $ws.Range("D15") '3.696.00'
Set-TextValue $ws.Range("E15") '  -1.03%  '

# Row 16
Set-TextValue $ws.Range("D16") '66.133.58'
Set-TextValue $ws.Range("E16") '  +0.10%  '

# Row 17
Set-TextValue $ws.Range("D17") '7.37'
Set-TextValue $ws.Range("E17") '  -1.68%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.179.31'
Set-TextValue $ws.Range("E18") '  -1.00%  '

# Row 19
Set-TextValue $ws.Range("E19") '  +1.06%  '

# Row 20
Set-TextValue $ws.Range("D20") '508.63'
Set-TextValue $ws.Range("E20") '  -0.53%  '

# Row 21
Set-TextValue $ws.Range("D21") '15.33'
Set-TextValue $ws.Range("E21") '  -0.82%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.727'
Set-TextValue $ws.Range("E22") '  -1.84%  '

# Row 23
Set-TextValue $ws.Range("D23") '7.97'
Set-TextValue $ws.Range("E23") '  -2.02%  '

# Row 24
Set-TextValue $ws.Range("D24") '14.74'
Set-TextValue $ws.Range("E24") '  -3.91%  '

# Row 25
Set-TextValue $ws.Range("D25") '84.39'
Set-TextValue $ws.Range("E25") '  -0.69%  '

# Row 26
Set-TextValue $ws.Range("E26") '  +0.05%  '

# Row 27
Set-TextValue $ws.Range("D27") '3.00'
Set-TextValue $ws.Range("E27") '  -0.61%  '

# Row 28
Set-TextValue $ws.Range("D28") '9.15'
Set-TextValue $ws.Range("E28") '  -1.69%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.37'
Set-TextValue $ws.Range("E29") '  +4.12%  '

# Row 30
Set-TextValue $ws.Range("D30") '3.00'
Set-TextValue $ws.Range("E30") '  +3.79%  '

# Row 31
Set-TextValue $ws.Range("D31") '7.18'
Set-TextValue $ws.Range("E31") '  +4.87%  '

# Row 32
Set-TextValue $ws.Range("D32") '27.90'
Set-TextValue $ws.Range("E32") '  -0.76%  '

# Row 33
Set-TextValue $ws.Range("E33") '  +0.20%  '

# Row 34
Set-TextValue $ws.Range("D34") '1.19'
Set-TextValue $ws.Range("E34") '  -3.00%  '

# Row 35
Set-TextValue $ws.Range("D35") '6.49'
Set-TextValue $ws.Range("E35") '  -1.39%  '

# Row 36
Set-TextValue $ws.Range("D36") '505.42'
Set-TextValue $ws.Range("E36") '  +3.96%  '

# Row 37
Set-TextValue $ws.Range("D37") '54.96'
Set-TextValue $ws.Range("E37") '  -0.08%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.0875'
Set-TextValue $ws.Range("E38") '  -3.44%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.0418'
Set-TextValue $ws.Range("E39") '  -0.55%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +5.51%  '

# Row 41
Set-TextValue $ws.Range("D41") '8.75'
Set-TextValue $ws.Range("E41") '  -1.97%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.0₃0682'
Set-TextValue $ws.Range("E42") '  +5.49%  '

# Row 43
Set-TextValue $ws.Range("D43") '2.84'
Set-TextValue $ws.Range("E43") '  -4.09%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.296'
Set-TextValue $ws.Range("E44") '  -1.78%  '

# Row 45
Set-TextValue $ws.Range("D45") '2.44'
Set-TextValue $ws.Range("E45") '  -0.40%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.826.20'
Set-TextValue $ws.Range("E46") '  -4.41%  '

# Row 47
Set-TextValue $ws.Range("D47") '27.92'
Set-TextValue $ws.Range("E47") '  -2.41%  '

# Row 48
Set-TextValue $ws.Range("B48") 'USDe'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D48") '0.999'
Set-TextValue $ws.Range("E48") '  -0.12%  '

# Row 49
Set-TextValue $ws.Range("B49") 'ThetaToken'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range("D49") '2.35'
Set-TextValue $ws.Range("E49") '  +1.45%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.116'
Set-TextValue $ws.Range("E50") '  +0.13%  '

# Row 51
Set-TextValue $ws.Range("B51") 'Arweave'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Range("D51") '35.20'
Set-TextValue $ws.Range("E51") '  +5.73%  '
